# Apply the edits described by the diff to the "python" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python")

# Row 64: keep existing Assignment(13) in B64, add "completed" to D64
$ws.Range("D64").Value = "completed"

# Row 65: add Assignment(13) to B65 and "completed" to D65
$ws.Range("B65").Value = "Assignment(13)"
$ws.Range("D65").Value = "completed"

# Row 66: add " Friday Holiday" (leading space) to B66
$ws.Range("B66").Value = " Friday Holiday"

# Row 67: change B67 from "Saturday" to "Assignment(13)", add "completed" to D67
$ws.Range("B67").Value = "Assignment(13)"
$ws.Range("D67").Value = "completed"

# Row 68 stays the same (A68=45725/Sunday date, B68=Sunday) - no change needed

# Row 69: new row - date 45726, "Python practice " (trailing space)
$ws.Range("A69").Value = 45726
$ws.Range("B69").Value = "Python practice "

# Row 70: new row - date 45727, "Assignment(5b)", completed
$ws.Range("A70").Value = 45727
$ws.Range("B70").Value = "Assignment(5b)"
$ws.Range("D70").Value = "completed"

# Row 71: new row - date 45728, "Assignment(5b)", completed
$ws.Range("A71").Value = 45728
$ws.Range("B71").Value = "Assignment(5b)"
$ws.Range("D71").Value = "completed"

# Copy the date style (numFmtId 14) from A68 onto the new date cells
$ws.Range("A68").Copy()
$ws.Range("A69:A71").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update sheet view to match target (top-left cell and selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 56
$ws.Range("A72").Select()
